$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sample file identifiers / cell ids for the two abuse-test rows
$ws.Range("A3").Value = "sample_2"
$ws.Range("B3").Value = "S2Abuse"

$ws.Range("B2").Value = "S1Abuse"
$ws.Range("A2").Value = "sample_1"

# Update the active selection shown when the sheet was last saved
$ws.Range("A9").Select()
